# build angular in DatingApp.API
# - Remove the "4. Implementacja" row (old row 9) from the Gantt table.
# - Rename "4.1 Implementacja wersji moblinej" -> "4. Implementacja aplikacji .Net Core"
# - Rename "4.2 Implementacja wersji WWW" -> "4.2 Implementacja Aplikacji Angular"
# - Update the start/end dates of several tasks.
# - Re-point the chart series to the now-shorter data range and shift the
#   chart up so it still starts right below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the row for "4. Implementacja" (row 9). Rows 10-13 shift up
#    to become rows 9-12.
# ------------------------------------------------------------------
$ws.Rows("9:9").Delete()

# ------------------------------------------------------------------
# 2. Rename the two sub-task labels that are now in rows 9 and 10.
# ------------------------------------------------------------------
$ws.Range("B9").Value = "4. Implementacja aplikacji .Net Core"
$ws.Range("B10").Value = "4.2 Implementacja Aplikacji Angular"

# ------------------------------------------------------------------
# 3. Update start (C) / end (D) dates for the affected rows.
#    (E column durations recalculate automatically.)
# ------------------------------------------------------------------
$ws.Range("C3").Value = 43736
$ws.Range("D3").Value = 43845

$ws.Range("C4").Value = 43736
$ws.Range("D4").Value = 43751

$ws.Range("C5").Value = 43763
$ws.Range("D5").Value = 43780

$ws.Range("C6").Value = 43763
$ws.Range("D6").Value = 43780

$ws.Range("C7").Value = 43765
$ws.Range("D7").Value = 43777

$ws.Range("C8").Value = 43766
$ws.Range("D8").Value = 43778

$ws.Range("D9").Value = 43849

$ws.Range("D10").Value = 43849

# ------------------------------------------------------------------
# 4. Re-point the chart's series formulas at the shrunk range
#    (was $B$3:$B$13 / $C$3:$C$13 / $E$3:$E$13, now ends at row 12).
# ------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(Sheet1!`$C`$2,Sheet1!`$B`$3:`$B`$12,Sheet1!`$C`$3:`$C`$12,1)"

$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = "=SERIES(Sheet1!`$D`$2,Sheet1!`$B`$3:`$B`$12,Sheet1!`$E`$3:`$E`$12,2)"

# Shift the chart up by the height of the deleted row so it keeps sitting
# right under the table (anchor row 13->12, 45->44).
$co.Top = $co.Top - $ws.Rows("9").Height

# ------------------------------------------------------------------
# 5. Update the active selection to match the saved state.
# ------------------------------------------------------------------
$ws.Range("D10").Select()
